# Limpieza de datos: corregidos errores de encoding y generados reportes
# Update counts and recomputed percentages in the "Datos_Faltantes" report
# (rows 2-11 of Sheet1), reflecting a denominator change from 6423 to 6420
# records after removing duplicate/bad rows from the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  B = 3925; C = 61.13707165109034 },
    @{ Row = 3;  B = 3807; C = 59.29906542056075 },
    @{ Row = 4;  B = 3214; C = 50.06230529595016 },
    @{ Row = 5;  B = 3197; C = 49.79750778816199 },
    @{ Row = 6;  B = 1939; C = 30.20249221183801 },
    @{ Row = 7;  B = 1939; C = 30.20249221183801 },
    @{ Row = 8;  B = 889;  C = 13.84735202492212 },
    @{ Row = 9;  B = 885;  C = 13.78504672897196 },
    @{ Row = 10; B = 13;   C = 0.2024922118380063 },
    @{ Row = 11; B = 13;   C = 0.2024922118380063 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
}
